$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "2026-02-19T17:00:13.963121+00:00"
$ws.Range("H2").Value = 33
$ws.Range("L2").Value = "[272285, 272284, 272296, 272297, 272295, 272303, 272301, 272305, 272302, 272312, 272313, 272308, 272315, 272316, 272304, 272399, 272407, 272403, 272405, 272408, 272413, 272490, 272501, 272500, 272592, 272600, 272605, 272601, 272610, 272673, 272667, 272674, 272675]"
